# Cyprus Division 1 - swap mis-ordered match rows (B:AD), keep column A (row index) fixed.
# Each group of rows below shares the same Date; the match data in columns B..AD
# was shifted by one position within the group. We restore the correct pairing
# by cyclically rotating B:AD among the rows in each group (row i <- row i+1,
# last row <- first row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    @(22, 23),
    @(50, 52),
    @(54, 55),
    @(67, 68),
    @(74, 75),
    @(111, 112),
    @(126, 127),
    @(138, 139),
    @(148, 149, 150),
    @(156, 158),
    @(169, 170),
    @(176, 177),
    @(223, 224),
    @(251, 252, 253),
    @(255, 257)
)

foreach ($rows in $groups) {
    $n = $rows.Count

    # Snapshot the current B:AD values for every row in this group first,
    # so later writes don't clobber data we still need to read.
    $snapshots = @()
    foreach ($r in $rows) {
        $snapshots += , ($ws.Range("B$r`:AD$r").Value2)
    }

    for ($i = 0; $i -lt $n; $i++) {
        $targetRow = $rows[$i]
        $srcValues = $snapshots[($i + 1) % $n]
        $ws.Range("B$targetRow`:AD$targetRow").Value2 = $srcValues
    }
}
